$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the GroundWater init scenario block from "18" references to "19"
# (row 5-11, column D holds the labels A18/B18/C18/G18/H18/I18/J18 -> A19/B19/C19/G19/H19/I19/J19)
$ws.Range("D5").Value  = "A19"
$ws.Range("D6").Value  = "B19"
$ws.Range("D7").Value  = "C19"
$ws.Range("D8").Value  = "G19"
$ws.Range("D9").Value  = "H19"
$ws.Range("D10").Value = "I19"
$ws.Range("D11").Value = "J19"

# Move the active selection to reflect where the editor left off
$ws.Range("D11").Select()
